# Weekly update: insert a new date's price records (Acelga, Vega Central
# Mapocho de Santiago) right after the previous week's block (old row 367),
# pushing all the existing records down by two rows. The two new records
# only have "Primera" and "Segunda" quality grades (no "Extra").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 368:369 - this shifts the old rows
# 368:482 down to 370:484 and extends the used range automatically.
$ws.Range("A368:R369").EntireRow.Insert()

# --- New row 368: Primera ---
$ws.Range("A368").Value = 9
$ws.Range("B368").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C368").Value = "Metropolitana"
$ws.Range("D368").Value = 44588
$ws.Range("E368").Value = 13
$ws.Range("F368").Value = 100112009
$ws.Range("G368").Value = "Acelga"
$ws.Range("H368").Value = "Sin especificar"
$ws.Range("I368").Value = "Primera"
$ws.Range("J368").Value = 61
$ws.Range("K368").Value = 15000
$ws.Range("L368").Value = 15000
$ws.Range("M368").Value = 15000
$ws.Range("N368").Value = "$/docena de atados"
$ws.Range("O368").Value = "Región Metropolitana"
$ws.Range("P368").Value = 5000
$ws.Range("Q368").Value = 3
$ws.Range("R368").Value = "Hortaliza"

# --- New row 369: Segunda ---
$ws.Range("A369").Value = 9
$ws.Range("B369").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C369").Value = "Metropolitana"
$ws.Range("D369").Value = 44588
$ws.Range("E369").Value = 13
$ws.Range("F369").Value = 100112009
$ws.Range("G369").Value = "Acelga"
$ws.Range("H369").Value = "Sin especificar"
$ws.Range("I369").Value = "Segunda"
$ws.Range("J369").Value = 43
$ws.Range("K369").Value = 13000
$ws.Range("L369").Value = 13000
$ws.Range("M369").Value = 13000
$ws.Range("N369").Value = "$/docena de atados"
$ws.Range("O369").Value = "Región Metropolitana"
$ws.Range("P369").Value = 4333
$ws.Range("Q369").Value = 3
$ws.Range("R369").Value = "Hortaliza"
